# edit.ps1
# Updates the "cryptos" worksheet data (rows 2-51) to reflect the latest
# scraped prices / percentages (and re-ordering of a couple of coins),
# per the commit "Updated cryptos list ... with GitHub Actions".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column keeps storing plain text (e.g. "31.006.91",
# "1.001", "0.000007746") instead of being auto-converted to numbers.
$ws.Range("D2:D51").NumberFormat = "@"

# Each entry: @(RankIndex, Coin, Link, Price, Volume(1h))
$data = @(
    @(0, 'Bitcoin', 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc', '31.006.91', '  +1.45%  '),
    @(1, 'Ethereum', 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth', '1.965.28', '  +2.61%  '),
    @(2, 'TetherUSD', 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt', '1.001', '  -0.05%  '),
    @(3, 'BNB', 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb', '247.89', '  +1.32%  '),
    @(4, 'USDC', 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc', '1.001', '  -0.05%  '),
    @(5, 'XRP', 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp', '0.4855', '  -0.16%  '),
    @(6, 'Cardano', 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada', '0.2957', '  +1.86%  '),
    @(7, 'Dogecoin', 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge', '0.06847', '  +1.66%  '),
    @(8, 'Solana', 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol', '19.18', '  -1.05%  '),
    @(9, 'Litecoin', 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc', '107.07', '  -3.50%  '),
    @(10, 'WrappedEther', 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth', '1.970.28', '  +2.78%  '),
    @(11, 'TRON', 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx', '0.07781', '  +2.76%  '),
    @(12, 'Polkadot', 'https://coinranking.com/coin/25W7FG7om+polkadot-dot', '5.452', '  +1.54%  '),
    @(13, 'Polygon', 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic', '0.7049', '  +4.93%  '),
    @(14, 'BitcoinCash', 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch', '286.96', '  -2.07%  '),
    @(15, 'WrappedBTC', 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc', '31.038.87', '  +1.52%  '),
    @(16, 'Avalanche', 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax', '13.25', '  +1.66%  '),
    @(17, 'ShibaInu', 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib', '0.000007746', '  +2.47%  '),
    @(18, 'WrappedliquidstakedEther2.0', 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth', '2.222.59', '  +2.19%  '),
    @(19, 'Uniswap', 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni', '5.634', '  +2.20%  '),
    @(20, 'Dai', 'https://coinranking.com/coin/MoTuySvg7+dai-dai', '1.002', '  +0.07%  '),
    @(21, 'BinanceUSD', 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd', '1.001', '  -0.08%  '),
    @(22, 'Chainlink', 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link', '6.615', '  +2.95%  '),
    @(23, 'Cosmos', 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom', '10.05', '  +6.02%  '),
    @(24, 'Monero', 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr', '170.80', '  +3.66%  '),
    @(25, 'EthereumClassic', 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc', '20.08', '  -1.39%  '),
    @(26, 'LidoDAOToken', 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo', '2.196', '  +4.44%  '),
    @(27, 'Stellar', 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm', '0.1068', '  -0.24%  '),
    @(28, 'Toncoin', 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton', '1.452', '  +1.25%  '),
    @(29, 'Filecoin', 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil', '4.845', '  +19.00%  '),
    @(30, 'InternetComputer(DFINITY)', 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp', '4.511', '  +8.97%  '),
    @(31, 'Hedera', 'https://coinranking.com/coin/jad286TjB+hedera-hbar', '0.05102', '  +1.86%  '),
    @(32, 'ImmutableX', 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx', '0.7750', '  +4.99%  '),
    @(33, 'ARBITRUM', 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb', '1.173', '  +3.07%  '),
    @(34, 'HuobiToken', 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht', '2.733', '  +0.61%  '),
    @(35, 'VeChain', 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet', '0.02043', '  +0.58%  '),
    @(36, 'MXToken', 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx', '2.727', '  +1.57%  '),
    @(37, 'FraxShare', 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs', '6.523', '  +11.38%  '),
    @(38, 'RenderToken', 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr', '2.132', '  +5.63%  '),
    @(39, 'TrustWalletToken', 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt', '0.8915', '  +3.20%  '),
    @(40, 'Quant', 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt', '110.01', '  +0.18%  '),
    @(41, 'TheSandbox', 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand', '0.4485', '  +0.80%  '),
    @(42, 'Aave', 'https://coinranking.com/coin/ixgUfzmLR+aave-aave', '72.72', '  +4.02%  '),
    @(43, 'PaxDollar', 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp', '1.001', '  -0.03%  '),
    @(44, 'Aptos', 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt', '7.540', '  +4.09%  '),
    @(45, 'Maker', 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr', '969.16', '  +14.87%  '),
    @(46, 'EnergySwap', 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens', '9.432', '  +1.74%  '),
    @(47, 'Algorand', 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo', '0.1268', '  +3.17%  '),
    @(48, 'Elrond', 'https://coinranking.com/coin/omwkOTglq+elrond-egld', '36.11', '  +3.22%  '),
    @(49, 'Decentraland', 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana', '0.4104', '  +1.57%  ')
)

$row = 2
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $ws.Cells.Item($row, 3).Value = $item[2]
    $ws.Cells.Item($row, 4).Value = $item[3]
    $ws.Cells.Item($row, 5).Value = $item[4]
    $row++
}

# Restore the default (unstyled) cell style now that the text values are
# safely stored, so the cell formatting matches the original workbook.
$ws.Range("D2:D51").Style = "Normal"
